$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column to remain text before writing number-like values,
# then reset the style back to Normal so no residual number-format style is left on cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.332.00"
$ws.Range("E2").Value = "  +2.36%  "

$ws.Range("D3").Value = "3.283.74"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "581.75"
$ws.Range("E5").Value = "  +4.58%  "

$ws.Range("D6").Value = "181.64"
$ws.Range("E6").Value = "  -1.81%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.275.33"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +1.97%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "0.574"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "46.14"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +5.20%  "

$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "636.20"
$ws.Range("E14").Value = "  +11.27%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.838.64"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "8.37"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.624.32"
$ws.Range("E17").Value = "  +2.93%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.118"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.308.08"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "17.54"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "10.85"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "0.891"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "17.62"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "5.01"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "97.39"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "3.97"
$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +4.27%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "32.55"
$ws.Range("E29").Value = "  +7.02%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "8.50"
$ws.Range("E30").Value = "  +0.60%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "6.65"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "590.26"
$ws.Range("E32").Value = "  +6.09%  "

$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "3.945.44"
$ws.Range("E33").Value = "  +4.89%  "

$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "3.58"
$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").Value = "10.87"
$ws.Range("E35").Value = "  +0.75%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.103"
$ws.Range("E36").Value = "  +1.05%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "55.60"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  +1.98%  "

$ws.Range("D40").Value = "3.23"
$ws.Range("E40").Value = "  +4.32%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").Value = "  +4.60%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "32.76"
$ws.Range("E42").Value = "  -2.09%  "

$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0683"
$ws.Range("E43").Value = "  +0.81%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "3.35"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  +1.36%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0413"
$ws.Range("E46").Value = "  +1.74%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.128"
$ws.Range("E47").Value = "  +1.43%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.53"
$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "1.32"
$ws.Range("E50").Value = "  +7.63%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "130.51"
$ws.Range("E51").Value = "  +4.81%  "

# Reset styles on the Price column back to Normal/default so text-coercion does not leave
# a stray number-format style applied to any cell.
$ws.Range("D2:D51").Style = "Normal"
